# "final stable before starting streamlit"
# - populate the "May 01" log sheet: last-updated stamp, balance formula,
#   consumed/burned totals, and a new "Bread" row
# - fan the day out into four more per-day sheets (May 011..May 014), each
#   seeded with just that same "Bread" row, mirroring the workflow that
#   will feed the upcoming Streamlit app

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("May 01")

# header / summary block
$ws.Range("B1").Value = "May 05 2022 20:08"
$ws.Range("B4").Value = 2200
$ws.Range("B5").Value = 2900
$ws.Range("B3").Formula = "=B4-B5"

# new food row - keep the "9.00"/"585.00"-style strings as text so the
# trailing zeroes survive (they'd collapse to plain numbers otherwise)
$ws.Range("B11:G11").NumberFormat = "@"
$ws.Range("A11").Value = "Bread"
$ws.Range("B11").Value = "9.00"
$ws.Range("C11").Value = "unit"
$ws.Range("D11").Value = "585.00"
$ws.Range("E11").Value = "37.80"
$ws.Range("F11").Value = "23.40"
$ws.Range("G11").Value = "30.60"

# four more day sheets, each just carrying the same Bread row - built by
# copying "May 01" (so they inherit its sheet formatting/margins) and then
# wiping rows 1-10 back out
$names = @("May 011", "May 012", "May 013", "May 014")
$prev = $ws
foreach ($n in $names) {
    $prev.Copy([System.Reflection.Missing]::Value, $prev) | Out-Null
    $newSheet = $wb.ActiveSheet
    $newSheet.Name = $n
    $newSheet.Rows("1:10").Clear() | Out-Null
    $newSheet.Range("A1").Select() | Out-Null
    $prev = $newSheet
}

$ws.Select() | Out-Null
